$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("C2").Value = 766207
$ws.Range("E2").Value = 1428874112
$ws.Range("C10").Value = 345453
$ws.Range("E10").Value = 1816729395
$ws.Range("C13").Value = 187763
$ws.Range("E13").Value = 1162522376
$ws.Range("C36").Value = 211191
$ws.Range("E36").Value = 404207239
$ws.Range("C57").Value = 31590
$ws.Range("E57").Value = 162173592
$ws.Range("C67").Value = 27092
$ws.Range("E67").Value = 168553916
$ws.Range("C72").Value = 331302
$ws.Range("E72").Value = 635354885
$ws.Range("C78").Value = 178406
$ws.Range("E78").Value = 892017619
$ws.Range("C79").Value = 680
$ws.Range("E79").Value = 20349120
$ws.Range("C91").Value = 18405
$ws.Range("E91").Value = 72117092
$ws.Range("C93").Value = 16555
$ws.Range("E93").Value = 48245887
$ws.Range("C112").Value = 145188
$ws.Range("E112").Value = 715584525
$ws.Range("C115").Value = 81784
$ws.Range("D115").Value = 14448
$ws.Range("E115").Value = 435971242
$ws.Range("C121").Value = 1305809
$ws.Range("E121").Value = 2273436091
$ws.Range("C127").Value = 9137
$ws.Range("E127").Value = 110235394
$ws.Range("C128").Value = 280
$ws.Range("E128").Value = 5719119
$ws.Range("C129").Value = 632760
$ws.Range("E129").Value = 3417328694
$ws.Range("C130").Value = 4227
$ws.Range("E130").Value = 139050789
$ws.Range("C132").Value = 585038
$ws.Range("E132").Value = 3441943057
$ws.Range("C136").Value = 26632
$ws.Range("E136").Value = 141844738
$ws.Range("C144").Value = 24513
$ws.Range("E144").Value = 88165149
$ws.Range("C151").Value = 39269
$ws.Range("E151").Value = 59810882
$ws.Range("C154").Value = 17975
$ws.Range("E154").Value = 69519959
$ws.Range("C157").Value = 630
$ws.Range("E157").Value = 1388087
$ws.Range("C171").Value = 95810
$ws.Range("E171").Value = 490330707
$ws.Range("C196").Value = 595464
$ws.Range("E196").Value = 983955743
$ws.Range("C215").Value = 230236
$ws.Range("E215").Value = 408662557
$ws.Range("C221").Value = 135467
$ws.Range("E221").Value = 681713588
$ws.Range("C229").Value = 612510
$ws.Range("E229").Value = 1040628728
$ws.Range("C237").Value = 283249
$ws.Range("E237").Value = 1437577049
$ws.Range("C240").Value = 205848
$ws.Range("E240").Value = 1066948938
